# Adds 8 new event rows (370-377) to the events table, mirroring the
# formatting of the template row above (row 369), and wires up the
# matching hyperlinks for the "Link" column (E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$templateRow = 369

# xlPasteFormats
$xlPasteFormats = -4122

function Set-PlainCell($row, $col, $value) {
    $ws.Cells.Item($templateRow, $col).Copy($ws.Cells.Item($row, $col))
    $ws.Cells.Item($row, $col).Value = $value
}

function Set-LinkCell($row, $url) {
    $cell = $ws.Cells.Item($row, 5)
    $cell.Value = $url
    $ws.Hyperlinks.Add($cell, $url, "", "", $url)
    # Re-apply the plain text/border formatting used by every other row in
    # the table (Hyperlinks.Add otherwise stamps its own "Hyperlink" style).
    $ws.Cells.Item($templateRow, 5).Copy()
    $cell.PasteSpecial($xlPasteFormats)
}

$rows = @(
    @{ Row = 370; Date = 45794; Event = "SUDHAUS RAVE"; Location = "Sudhaus"; Stadt = "Unna"; Link = "https://www.instagram.com/reel/DIgw7tgMx6T/?igsh=MWxlN3o3ZmFkZHQxMA==" },
    @{ Row = 371; Date = 45786; Event = "FLYM x H33 SHOWCASE"; Location = "SNRS"; Stadt = "Dortmund"; Link = "https://www.instagram.com/p/DIrShFeqeEO/?igsh=MW92eDFxdXdzd2ZzZA==" },
    @{ Row = 372; Date = 45772; Event = "SOLI RAVE"; Location = "Die Nacht"; Stadt = "Mönchengladbach"; Link = "https://www.instagram.com/reel/DIwfEBkNX69/?igsh=aTZoNzk1ZzZvZTNx" },
    @{ Row = 373; Date = 45793; Event = "SUBSURFACE"; Location = "Stollen134"; Stadt = "Dortmund"; Link = "https://www.instagram.com/reel/DIwVRsGoa2e/?igsh=MWFpMDNmMzEwZWw1NQ==" },
    @{ Row = 374; Date = 45794; Event = "FEIERN FÜR EINEN GUTEN ZWECK"; Location = "Samy's"; Stadt = "Düsseldorf"; Link = "https://www.instagram.com/p/DIwETyHAz_3/?igsh=MXFrMTd1b2ZyMDdxbA==" },
    @{ Row = 375; Date = 45785; Event = "TURBO 120 MIN RAVE"; Location = "Oma Doris"; Stadt = "Dortmund"; Link = "https://www.instagram.com/reel/DIwcQSNKgda/?igsh=NDF6a3Nkb2pjcjFx" },
    @{ Row = 376; Date = 45773; Event = "POLLERWIESEN OPENING NIGHT x BOOTSHAUS"; Location = "Bootshaus"; Stadt = "Köln"; Link = "https://bootshaus.tv/tickets/" },
    @{ Row = 377; Date = 45815; Event = "SMILEYVENTS ROOFTOP EDITION"; Location = "check event link"; Stadt = "Duisburg"; Link = "https://www.instagram.com/reel/DIuDJ-bOV7a/?igsh=YnE2NnU3Y2YxZ3hw" }
)

foreach ($r in $rows) {
    Set-PlainCell $r.Row 1 $r.Date
    Set-PlainCell $r.Row 2 $r.Event
    Set-PlainCell $r.Row 3 $r.Location
    Set-PlainCell $r.Row 4 $r.Stadt
    Set-LinkCell $r.Row $r.Link
}
